# Merge Redhat and Ubuntu smoketesting rows into a single "Linux" row
# on the "OS instructions" sheet (issue_template.xlsx, PR #107).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OS instructions")

# Row 4 used to be "Redhat" - it becomes the merged "Linux" row, reusing
# the additional-instructions text that used to live on the "Ubuntu" row
# (the tar.xz / conda install instructions for Linux).
$ws.Cells.Item(4, 1).Value = "Linux"
$ws.Cells.Item(4, 2).Value = "* To install the tar.xz package for Linux, run ``(sudo) tar -xJf mantid-VA.B.C.tar.xz`` in a terminal and it will unzip the package in your current working directory. 
* To install via conda:
  - Use Intel Conda and make sure conda-forge is added to channels
  - In terminal, create a new empty environment and activate it
  - run ``conda install -c ""mantid/label/vA.B.C-rc1"" mantidworkbench`` , where A.B.C is the release version.
"

# Row 5 used to be "Ubuntu" - now that it has been folded into row 4 it is
# left blank.
$ws.Cells.Item(5, 1).Value = ""
$ws.Cells.Item(5, 2).Value = ""

# Row heights re-wrap to fit the new (now longer/shorter) text.
$ws.Rows.Item(4).RowHeight = 132
$ws.Rows.Item(5).RowHeight = 20

# Selection moved as part of the author's edit.
$ws.Range("B7").Select()
